# Refresh the scraped crypto Price (D) / 1h Volume change (E) columns,
# row 2 through row 51, with the newly polled values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.344.12"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "3.122.19"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'567.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "'147.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "3.123.16"
$ws.Range("E8").Value = "  +0.62%  "
$ws.Range("D9").Value = "'0.520"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").Value = "'0.156"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.13%  "
$ws.Range("D11").Value = "'6.02"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.23%  "
$ws.Range("D12").Value = "'0.489"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").Value = "'0.0000257"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.44%  "
$ws.Range("D14").Value = "'36.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").Value = "3.633.19"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("D16").Value = "64.619.58"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").Value = "3.137.24"
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("D18").Value = "'6.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").Value = "'0.110"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").Value = "'495.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").Value = "'14.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.74%  "
$ws.Range("D22").Value = "'0.703"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("D23").Value = "'14.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.54%  "
$ws.Range("D24").Value = "'7.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.84%  "
$ws.Range("D25").Value = "'83.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("D26").Value = "'0.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("D27").Value = "'2.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("D28").Value = "'8.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.82%  "
$ws.Range("D29").Value = "'2.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").Value = "'27.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.17%  "
$ws.Range("D31").Value = "'2.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.41%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("D33").Value = "'1.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").Value = "'6.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.23%  "
$ws.Range("D35").Value = "'6.36"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.61%  "
$ws.Range("D36").Value = "'54.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("D37").Value = "'0.0885"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.23%  "
$ws.Range("D38").Value = "'459.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "'0.0411"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("D40").Value = "'2.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.23%  "
$ws.Range("D41").Value = "'8.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").Value = "3.007.08"
$ws.Range("E42").Value = "  -2.73%  "
$ws.Range("D43").Value = "'0.115"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.93%  "
$ws.Range("D44").Value = "'0.278"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.62%  "
$ws.Range("D45").Value = "'2.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.42%  "
$ws.Range("D46").Value = "'27.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.63%  "
$ws.Range("D47").Value = "0.0₃0565"
$ws.Range("E47").Value = "  +4.55%  "
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("D49").Value = "'0.113"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.80%  "
$ws.Range("D50").Value = "'2.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("D51").Value = "'117.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.13%  "
